$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# -----------------------------------------------------------------
# Edit 1: {sutotaNet}€  ->  {subtotalNet}€   (row 6, col 4)
#   Split into 5 runs with identical rPr: "{su" / "b" / "tota" / "l" / "Net}€"
# -----------------------------------------------------------------
$cell1 = $t.Cell(6, 4)
$range1 = $cell1.Range
$start1 = $range1.Start

# Insert the two new characters first (plain text, no formatting yet)
# Insert 'l' before 'N' (offset 7 in "{sutotaNet}€")
$d.Range($start1 + 7, $start1 + 7).InsertAfter("l")
# Insert 'b' after "{su" (offset 3)
$d.Range($start1 + 3, $start1 + 3).InsertAfter("b")

# Now the text reads "{subtotalNet}€". Force run boundaries around each
# inserted character by toggling Bold off/on (it was already Bold, so the
# effective formatting is unchanged, but a new run boundary is created).
$bRange = $d.Range($start1 + 3, $start1 + 4)
$bRange.Font.Bold = 0
$bRange.Font.Bold = 1

$lRange = $d.Range($start1 + 8, $start1 + 9)
$lRange.Font.Bold = 0
$lRange.Font.Bold = 1

# -----------------------------------------------------------------
# Edit 2: {totalCrossPrice}€  ->  {totalGrossPrice}€   (row 9, col 4)
#   Split into 3 runs with identical rPr: "{total" / "G" / "rossPrice}€"
#   (the leading space run before the placeholder must stay untouched/separate)
# -----------------------------------------------------------------
$cell2 = $t.Cell(9, 4)
$range2 = $cell2.Range
$start2 = $range2.Start

# Replace the 'C' (offset 7 in " {totalCrossPrice}€") with 'G'
$d.Range($start2 + 7, $start2 + 8).Text = "G"

# Force a run boundary around the replaced character
$gRange = $d.Range($start2 + 7, $start2 + 8)
$gRange.Font.Bold = 1
$gRange.Font.Bold = 0

# Keep the leading space run separate from the (re-split) placeholder run
$spaceRange = $d.Range($start2, $start2 + 1)
$spaceRange.Font.Bold = 1
$spaceRange.Font.Bold = 0
